$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$tblShape = $s1.Shapes.Item(1)
$tbl = $tblShape.Table
$cell = $tbl.Cell(2, 1)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Font.Name = $tr.Font.Name
